# Updated capital structure database
# Applies updated values to rows 2 and 3 of the "earnings_debt" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$values = @{
    "D"  = -0.452
    "G"  = -37.11111111111111
    "H"  = -37.11111111111111
    "I"  = -50.33333333333334
    "J"  = -50.33333333333334
    "K"  = -0.8
    "L"  = -88.8888888888889
    "V"  = 0.0004484304932735426
    "W"  = -0.3902439024390245
    "X"  = 0.3078454964923482
    "Y"  = -0.6980893989313726
    "Z"  = 0.0008629782337712148
    "AA" = -0.04343657109981781
    "AB" = 0.1089274929401542
    "AC" = -0.152364064039972
    "AD" = 8.23
    "AF" = 8.23
    "AG" = 8.229000000000001
    "AH" = 0.7868068833652008
    "AI" = 0.8811563169164883
    "AJ" = 0.78678649966536
    "AK" = 0.8811435913909412
    "AL" = 0.302
    "AM" = 0.302
    "AN" = -25.09146341463414
    "AO" = -1.5
    "AP" = -25.08841463414634
    "AQ" = -1.5
}

foreach ($row in 2,3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
